$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 11906440
$ws.Range("I100").Value = 15873910
$ws.Range("J100").Value = 4030.2856
$ws.Range("K100").Value = 15873910
$ws.Range("L100").Value = 4030.2856
$ws.Range("M100").Value = -15873369
$ws.Range("N100").Value = -5112.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 3143.818
$ws.Range("J27").Value = 3143.818
$ws.Range("L27").Value = 3143.818
$ws.Range("N27").Value = -3511.818
$ws.Range("H32").Value = 492413.56
$ws.Range("I32").Value = 5191.884
$ws.Range("J32").Value = 1801821.9
$ws.Range("K32").Value = 5191.884
$ws.Range("L32").Value = 1801821.9
$ws.Range("M32").Value = -4904.884
$ws.Range("N32").Value = -1802395.9
$ws.Range("H61").Value = 2357.25
$ws.Range("I61").Value = 2455.8948
$ws.Range("K61").Value = 2455.8948
$ws.Range("M61").Value = -2243.8948
$ws.Range("H97").Value = 1481.8182
$ws.Range("I97").Value = 1050
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1050
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -554
$ws.Range("N97").Value = -2992
$ws.Range("H136").Value = 2357.25
$ws.Range("I136").Value = 2455.8948
$ws.Range("K136").Value = 7367.6844
$ws.Range("M136").Value = -4817.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2334.7144
$ws.Range("I99").Value = 1242.2941
$ws.Range("J99").Value = 4023
$ws.Range("K99").Value = 1242.2941
$ws.Range("L99").Value = 4023
$ws.Range("M99").Value = 255.7058999999999
$ws.Range("N99").Value = -7019
$ws.Range("H134").Value = 1869.2041
$ws.Range("I134").Value = 1783.186
$ws.Range("J134").Value = 2485.6667
$ws.Range("K134").Value = 5349.558
$ws.Range("L134").Value = 7457.000100000001
$ws.Range("M134").Value = -2814.558
$ws.Range("N134").Value = -12527.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1845.7587
$ws.Range("I31").Value = 2163.4707
$ws.Range("J31").Value = 1395.6666
$ws.Range("K31").Value = 2163.4707
$ws.Range("L31").Value = 1395.6666
$ws.Range("M31").Value = -1868.4707
$ws.Range("N31").Value = -1985.6666
$ws.Range("H34").Value = 1845.7587
$ws.Range("I34").Value = 2163.4707
$ws.Range("J34").Value = 1395.6666
$ws.Range("K34").Value = 2163.4707
$ws.Range("L34").Value = 1395.6666
$ws.Range("M34").Value = -1961.4707
$ws.Range("N34").Value = -1799.6666
$ws.Range("H58").Value = 40001316
$ws.Range("I58").Value = 76923920
$ws.Range("J58").Value = 1830.3334
$ws.Range("K58").Value = 76923920
$ws.Range("L58").Value = 1830.3334
$ws.Range("M58").Value = -76923717
$ws.Range("N58").Value = -2236.3334
$ws.Range("H62").Value = 2611.25
$ws.Range("I62").Value = 2566.6667
$ws.Range("J62").Value = 2745
$ws.Range("K62").Value = 2566.6667
$ws.Range("L62").Value = 2745
$ws.Range("M62").Value = -1942.6667
$ws.Range("N62").Value = -3993
$ws.Range("H65").Value = 2611.25
$ws.Range("I65").Value = 2566.6667
$ws.Range("J65").Value = 2745
$ws.Range("K65").Value = 12833.3335
$ws.Range("L65").Value = 13725
$ws.Range("M65").Value = -9713.333500000001
$ws.Range("N65").Value = -19965
$ws.Range("H132").Value = 6252100.5
$ws.Range("I132").Value = 1596.5
$ws.Range("K132").Value = 4789.5
$ws.Range("M132").Value = -2259.5
$ws.Range("H134").Value = 50007256
$ws.Range("I134").Value = 75000960
$ws.Range("K134").Value = 225002880
$ws.Range("M134").Value = -225000345
$ws.Range("H136").Value = 40001316
$ws.Range("I136").Value = 76923920
$ws.Range("J136").Value = 1830.3334
$ws.Range("K136").Value = 230771760
$ws.Range("L136").Value = 5491.0002
$ws.Range("M136").Value = -230769210
$ws.Range("N136").Value = -10591.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 875.1698
$ws.Range("I68").Value = 649.75
$ws.Range("J68").Value = 915.24445
$ws.Range("K68").Value = 1949.25
$ws.Range("L68").Value = 2745.73335
$ws.Range("M68").Value = -1138.25
$ws.Range("N68").Value = -4367.73335
$ws.Range("H71").Value = 875.1698
$ws.Range("I71").Value = 649.75
$ws.Range("J71").Value = 915.24445
$ws.Range("K71").Value = 5847.75
$ws.Range("L71").Value = 8237.200049999999
$ws.Range("M71").Value = -1791.75
$ws.Range("N71").Value = -16349.20005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2000802.2
$ws.Range("I21").Value = 2000802.2
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 2000802.2
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = -2000629.2
$ws.Range("H30").Value = 2000802.2
$ws.Range("I30").Value = 2000802.2
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2000802.2
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = -2000697.2
$ws.Range("H43").Value = 76600
$ws.Range("J43").Value = 76600
$ws.Range("L43").Value = 76600
$ws.Range("N43").Value = -76902
$ws.Range("H102").Value = 2137.4285
$ws.Range("I102").Value = 2210.3333
$ws.Range("J102").Value = 1700
$ws.Range("K102").Value = 2210.3333
$ws.Range("L102").Value = 1700
$ws.Range("M102").Value = -588.3332999999998
$ws.Range("N102").Value = -4944

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 501.68
$ws.Range("I22").Value = 450.55884
$ws.Range("J22").Value = 610.3125
$ws.Range("K22").Value = 450.55884
$ws.Range("L22").Value = 610.3125
$ws.Range("M22").Value = -155.55884
$ws.Range("N22").Value = -1200.3125
$ws.Range("H27").Value = 501.68
$ws.Range("I27").Value = 450.55884
$ws.Range("J27").Value = 610.3125
$ws.Range("K27").Value = 450.55884
$ws.Range("L27").Value = 610.3125
$ws.Range("M27").Value = -343.55884
$ws.Range("N27").Value = -824.3125
$ws.Range("H104").Value = 15498.823
$ws.Range("J104").Value = 15498.823
$ws.Range("L104").Value = 15498.823
$ws.Range("N104").Value = -22486.823
$ws.Range("H106").Value = 20370
$ws.Range("J106").Value = 20370
$ws.Range("L106").Value = 20370
$ws.Range("N106").Value = -22894
$ws.Range("H132").Value = 3334992
$ws.Range("I132").Value = 4762917.5
$ws.Range("K132").Value = 14288752.5
$ws.Range("M132").Value = -14286222.5
$ws.Range("H136").Value = 38444.605
$ws.Range("I136").Value = 83951.336
$ws.Range("J136").Value = 4314.5625
$ws.Range("K136").Value = 251854.008
$ws.Range("L136").Value = 12943.6875
$ws.Range("M136").Value = -249304.008
$ws.Range("N136").Value = -18043.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3166.6667
$ws.Range("I62").Value = 2875
$ws.Range("J62").Value = 3750
$ws.Range("K62").Value = 2875
$ws.Range("L62").Value = 3750
$ws.Range("M62").Value = -2251
$ws.Range("N62").Value = -4998
$ws.Range("H65").Value = 3166.6667
$ws.Range("I65").Value = 2875
$ws.Range("J65").Value = 3750
$ws.Range("K65").Value = 14375
$ws.Range("L65").Value = 18750
$ws.Range("M65").Value = -11255
$ws.Range("N65").Value = -24990
$ws.Range("H104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("L104").Value = 10000
$ws.Range("N104").Value = -16988
$ws.Range("H132").Value = 50002180
$ws.Range("I132").Value = 64518080
$ws.Range("J132").Value = 2972
$ws.Range("K132").Value = 193554240
$ws.Range("L132").Value = 8916
$ws.Range("M132").Value = -193551710
$ws.Range("N132").Value = -13976
$ws.Range("H136").Value = 818.0417
$ws.Range("I136").Value = 329
$ws.Range("K136").Value = 987
$ws.Range("M136").Value = 1563
